$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows starting at row 17 so the ILO source block
# (currently rows 20-21) shifts down to rows 26-27.
$ws.Range("A17:A22").EntireRow.Insert()

# New "MSME definition" header row (row 17) - bold "title" style like other headers
$ws.Range("B17").Value = "Number of employees"
$ws.Range("C17").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D17").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B17:D17").Font.Bold = $true

# Data rows 18-21 (normal style)
$ws.Range("A18").Value = "Micro"
$ws.Range("B18").Value = "<5"
$ws.Range("C18").Formula = "="""""
$ws.Range("D18").Formula = "="""""

$ws.Range("A19").Value = "Small"
$ws.Range("B19").Value = "5-9"
$ws.Range("C19").Formula = "="""""
$ws.Range("D19").Formula = "="""""

$ws.Range("A20").Value = "Medium"
$ws.Range("B20").Value = "10-29"
$ws.Range("C20").Formula = "="""""
$ws.Range("D20").Formula = "="""""

$ws.Range("A21").Value = "Large"
$ws.Range("B21").Value = ">=30"
$ws.Range("C21").Formula = "="""""
$ws.Range("D21").Formula = "="""""

Write-Host "done"
